$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 129, shifting rows 129:241 down to 130:242.
$ws.Rows(129).Insert()

# Populate the new row 129 with the new record.
$ws.Range("A129").Value = 7
$ws.Range("B129").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C129").Value = "Ñuble"
$ws.Range("D129").Value = 45040
$ws.Range("D129").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E129").Value = 16
$ws.Range("F129").Value = "Fruta"
$ws.Range("G129").Value = 100101
$ws.Range("H129").Value = "Berries"
$ws.Range("I129").Value = 100101007
$ws.Range("J129").Value = "Kiwi"
$ws.Range("K129").Value = "Hayward"
$ws.Range("L129").Value = "Primera"
$ws.Range("M129").Value = 80
$ws.Range("N129").Value = 13000
$ws.Range("O129").Value = 13000
$ws.Range("P129").Value = 13000
$ws.Range("Q129").Value = "$/bandeja 18 kilos"
$ws.Range("R129").Value = "Región de O'Higgins"
$ws.Range("S129").Value = 722
$ws.Range("T129").Value = 18
